$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

$ws.Range("C3").Value  = "легк сер б/к"
$ws.Range("C4").Value  = "легк сер б/к"
$ws.Range("C10").Value = "210B C H сер Type"
$ws.Range("C11").Value = "LS-2 груз сер Type"
$ws.Range("C12").Value = "202B C сер Type"
$ws.Range("C13").Value = "202B C LS-2 H сер Type"
$ws.Range("C14").Value = "груз сер б/к"
$ws.Range("C15").Value = "легк сер б/к"
$ws.Range("C16").Value = "легк сер б/к"
